# Actualización 10 de Mayo
# Updates statistics on the three "Estadisticos" sheets (1P, 2P, Final) to
# reflect a student (ZACARIAS HERNANDEZ LUIS ALBERTO, grupo 6ARHV) moving
# from "Blancos" (pending) into the graded pool, and registers that
# student on the "Rescatables" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Estadisticos 1P": columns D=Blancos, E=Reprobados, F=Aprobados,
# G=Por_Apro, H=Promedio
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Estadisticos 1P")

$ws1.Range("D2").Value = 0
$ws1.Range("E2").Value = 14
$ws1.Range("F2").Value = 15
$ws1.Range("G2").Value = 51.72
$ws1.Range("H2").Value = 6.7

$ws1.Range("D3").Value = 0
$ws1.Range("E3").Value = 11
$ws1.Range("F3").Value = 10
$ws1.Range("G3").Value = 47.62
$ws1.Range("H3").Value = 6.2

$ws1.Range("D4").Value = 0
$ws1.Range("E4").Value = 7
$ws1.Range("F4").Value = 21
$ws1.Range("G4").Value = 75
$ws1.Range("H4").Value = 7.3

$ws1.Range("D5").Value = 0
$ws1.Range("E5").Value = 15
$ws1.Range("H5").Value = 5.7

$ws1.Range("D6").Value = 0
$ws1.Range("E6").Value = 15
$ws1.Range("F6").Value = 7
$ws1.Range("G6").Value = 31.82
$ws1.Range("H6").Value = 5.5

# ---------------------------------------------------------------------
# Sheet "Estadisticos 2P": same column layout; H (Promedio) was blank
# before and now gets values too.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Estadisticos 2P")

$ws2.Range("D2").Value = 0
$ws2.Range("F2").Value = 16
$ws2.Range("G2").Value = 55.17
$ws2.Range("H2").Value = 6.7

$ws2.Range("D3").Value = 0
$ws2.Range("E3").Value = 11
$ws2.Range("F3").Value = 10
$ws2.Range("G3").Value = 47.62
$ws2.Range("H3").Value = 6.2

$ws2.Range("D4").Value = 0
$ws2.Range("E4").Value = 5
$ws2.Range("F4").Value = 23
$ws2.Range("G4").Value = 82.14
$ws2.Range("H4").Value = 7.3

$ws2.Range("D5").Value = 0
$ws2.Range("E5").Value = 14
$ws2.Range("F5").Value = 8
$ws2.Range("G5").Value = 36.36
$ws2.Range("H5").Value = 5.7

$ws2.Range("D6").Value = 0
$ws2.Range("E6").Value = 15
$ws2.Range("F6").Value = 7
$ws2.Range("G6").Value = 31.82
$ws2.Range("H6").Value = 5.5

# ---------------------------------------------------------------------
# Sheet "Estadisticos Final": same column layout.
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Estadisticos Final")

$ws3.Range("D2").Value = 0
$ws3.Range("E2").Value = 13
$ws3.Range("F2").Value = 16
$ws3.Range("G2").Value = 55.17
$ws3.Range("H2").Value = 6.8

$ws3.Range("D3").Value = 0
$ws3.Range("E3").Value = 11
$ws3.Range("F3").Value = 10
$ws3.Range("G3").Value = 47.62
$ws3.Range("H3").Value = 6

$ws3.Range("D4").Value = 0
$ws3.Range("E4").Value = 5
$ws3.Range("F4").Value = 23
$ws3.Range("G4").Value = 82.14
$ws3.Range("H4").Value = 7.4

$ws3.Range("D5").Value = 0
$ws3.Range("E5").Value = 14
$ws3.Range("F5").Value = 8
$ws3.Range("G5").Value = 36.36
$ws3.Range("H5").Value = 6

$ws3.Range("D6").Value = 0
$ws3.Range("E6").Value = 15
$ws3.Range("F6").Value = 7
$ws3.Range("G6").Value = 31.82
$ws3.Range("H6").Value = 5.8

# ---------------------------------------------------------------------
# Sheet "Rescatables": add the new student row.
# Columns: A=NC, B=Paterno, C=Materno, D=Nombres, E=Nombre_Largo,
# F=Grupo, G=Reprobadas
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Rescatables")

$ws4.Range("A2").Value = 18330051920391
$ws4.Range("B2").Value = "ZACARIAS"
$ws4.Range("C2").Value = "HERNANDEZ"
$ws4.Range("D2").Value = "LUIS ALBERTO"
$ws4.Range("E2").Value = "PROBABILIDAD Y ESTADÍSTICA"
$ws4.Range("F2").Value = "6ARHV"
$ws4.Range("G2").Value = 2
